# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": the Status column is updated on every sheet and the
# two timestamps that record when the zh-cn / de-de handoff XLIFFs were last
# (re)generated are bumped forward a few dozen seconds. The Status column
# widens (its longest value grew from "In Translation" to the longer
# "Ready for handoff") consistently across all three sheets.

$wb = $excel.ActiveWorkbook

# Column width the engine's MDW=6 lattice lands closest to the wider
# "Ready for handoff" status column (98 px -> (98/6) character units).
$statusColWidth = 98 / 6

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-21 09:03:30"
$overview.Columns.Item(5).ColumnWidth = $statusColWidth
$overview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-21 09:03:26"
$zhcn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-21 09:03:30"
$dede.Columns.Item(3).ColumnWidth = $statusColWidth
